# Guards can now be dragged along with the Player
# Adds a new row (row 6) to Tabelle1 documenting that guards can be
# dragged/pushed, mirroring the formatting of the existing date/notes rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6, column A: new test date (2017-09-07 == serial 42985), formatted
# the same way as the existing date cells (A2/A5) - copy format first so
# we reuse the existing date/alignment style instead of creating a new one.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 42985

# Row 6, column D: the new "Zu Verbessern" note - copy formatting from the
# cell above (D5) so it reuses the existing wrap-text style.
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = "Wachen können gezogen/geschoben werden"

# Match the row height used by the other wrapped-text note rows.
$ws.Rows.Item(6).RowHeight = 30

# Move the active selection down, like after typing the new row in Excel.
[void]$ws.Range("A7").Select()
